$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Full-Delivery": update a couple of values and move the selection.
# ---------------------------------------------------------------------------
$wsFull = $wb.Worksheets.Item("Full-Delivery")
$wsFull.Range("B49").Value = 19
$wsFull.Range("E50").Value = 65
[void]$wsFull.Range("E51").Select()

# ---------------------------------------------------------------------------
# Sheet "Progress": update row 3, 5, 6, 7 values and move the selection.
# ---------------------------------------------------------------------------
$wsProgress = $wb.Worksheets.Item("Progress")
$wsProgress.Range("G3").Value = 19
$wsProgress.Range("G5:AR5").Value = 74
$wsProgress.Range("G6").Value = -11
$wsProgress.Range("G7:AR7").Value = 695
[void]$wsProgress.Range("E13").Select()

# ---------------------------------------------------------------------------
# Sheet "Sprint 31": update row 24 values.
# ---------------------------------------------------------------------------
$wsSprint31 = $wb.Worksheets.Item("Sprint 31")
$wsSprint31.Range("B24").Value = 16.25
$wsSprint31.Range("C24").Value = 16.25
$wsSprint31.Range("F24").Value = 16.25
$wsSprint31.Range("G24").Value = 16.25

# ---------------------------------------------------------------------------
# Restore "Full-Delivery" as the active sheet (it was active before the
# edits and the diff does not indicate the active tab changed).
# ---------------------------------------------------------------------------
[void]$wsFull.Activate()
